$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.709.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.060.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.666'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.35'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '60.93'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.369'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0752'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.88%  '
$ws.Range("E12").Value = '  -3.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.939'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.87'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.363.23'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.32%  '
$ws.Range("E16").Value = '  -3.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.050.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.03%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.28%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.619.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.97%  '
$ws.Range("E20").Value = '  -2.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0865'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.89'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.63%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  -2.47%  '
$ws.Range("E26").Value = '  +6.07%  '
$ws.Range("E27").Value = '  -6.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '166.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.51'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.09%  '
$ws.Range("E34").Value = '  -2.82%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.84'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0848'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.22'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.08'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.63%  '
$ws.Range("E40").Value = '  -5.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.19%  '
$ws.Range("E42").Value = '  -3.02%  '
$ws.Range("E43").Value = '  -4.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '95.16'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0910'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.416.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +12.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.84%  '
$ws.Range("E49").Value = '  +1.63%  '
$ws.Range("E50").Value = '  -3.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.249.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.31%  '
